$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.198.09"
$ws.Range("E2").Value = "  +5.42%  "

$ws.Range("D3").Value = "1.784.45"
$ws.Range("E3").Value = "  +3.05%  "

$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").Value = "'244.57"
$ws.Range("E5").Value = "  +0.85%  "

$ws.Range("D7").Value = "'0.4917"
$ws.Range("E7").Value = "  -0.30%  "

$ws.Range("D8").Value = "'0.2680"
$ws.Range("E8").Value = "  +2.08%  "

$ws.Range("D9").Value = "'0.06273"
$ws.Range("E9").Value = "  +0.83%  "

$ws.Range("D10").Value = "1.779.35"
$ws.Range("E10").Value = "  +2.76%  "

$ws.Range("D11").Value = "'16.50"
$ws.Range("E11").Value = "  +3.29%  "

$ws.Range("D12").Value = "'0.07041"
$ws.Range("E12").Value = "  +0.65%  "

$ws.Range("D13").Value = "'0.6283"
$ws.Range("E13").Value = "  +2.56%  "

$ws.Range("E14").Value = "  +3.29%  "

$ws.Range("D15").Value = "'79.98"
$ws.Range("E15").Value = "  +3.45%  "

$ws.Range("D16").Value = "28.167.64"
$ws.Range("E16").Value = "  +6.16%  "

$ws.Range("E17").Value = "  +0.24%  "

$ws.Range("D18").Value = "'0.9994"
$ws.Range("E18").Value = "  +0.15%  "

$ws.Range("D19").Value = "'0.000007249"
$ws.Range("E19").Value = "  +0.26%  "

$ws.Range("E20").Value = "  +5.15%  "

$ws.Range("D21").Value = "2.007.17"
$ws.Range("E21").Value = "  +2.97%  "

$ws.Range("D22").Value = "'4.553"
$ws.Range("E22").Value = "  +1.27%  "

$ws.Range("D23").Value = "'8.753"
$ws.Range("E23").Value = "  +2.00%  "

$ws.Range("D24").Value = "'5.268"
$ws.Range("E24").Value = "  +3.10%  "

$ws.Range("D25").Value = "'141.32"
$ws.Range("E25").Value = "  +2.13%  "

$ws.Range("E26").Value = "  +2.80%  "

$ws.Range("E27").Value = "  +4.60%  "

$ws.Range("D28").Value = "'109.70"
$ws.Range("E28").Value = "  +2.92%  "

$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("D30").Value = "'4.186"
$ws.Range("E30").Value = "  +6.25%  "

$ws.Range("D31").Value = "'0.08274"
$ws.Range("E31").Value = "  +3.50%  "

$ws.Range("D32").Value = "'3.763"
$ws.Range("E32").Value = "  +2.27%  "

$ws.Range("D33").Value = "'0.04895"
$ws.Range("E33").Value = "  +9.16%  "

$ws.Range("D34").Value = "'1.079"
$ws.Range("E34").Value = "  +7.50%  "

$ws.Range("D35").Value = "'2.615"
$ws.Range("E35").Value = "  +0.22%  "

$ws.Range("D36").Value = "'0.6521"
$ws.Range("E36").Value = "  +4.36%  "

$ws.Range("D37").Value = "'0.9502"
$ws.Range("E37").Value = "  +1.03%  "

$ws.Range("D38").Value = "'2.611"
$ws.Range("E38").Value = "  +7.73%  "

$ws.Range("D39").Value = "'2.047"
$ws.Range("E39").Value = "  -0.20%  "

$ws.Range("D40").Value = "'5.921"
$ws.Range("E40").Value = "  +5.93%  "

$ws.Range("D41").Value = "'0.01554"
$ws.Range("E41").Value = "  +2.56%  "

$ws.Range("E42").Value = "  +0.15%  "

$ws.Range("D43").Value = "'99.99"
$ws.Range("E43").Value = "  +0.56%  "

$ws.Range("D44").Value = "'0.3985"
$ws.Range("E44").Value = "  +3.05%  "

$ws.Range("D45").Value = "'7.187"
$ws.Range("E45").Value = "  +3.24%  "

$ws.Range("D46").Value = "'0.1217"
$ws.Range("E46").Value = "  +4.85%  "

$ws.Range("D47").Value = "'0.05446"
$ws.Range("E47").Value = "  +1.14%  "

$ws.Range("D48").Value = "'8.028"
$ws.Range("E48").Value = "  +1.47%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'30.77"
$ws.Range("E49").Value = "  +1.75%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.294"
$ws.Range("E50").Value = "  +4.55%  "

$ws.Range("E51").Value = "  +2.05%  "

